$d = $word.ActiveDocument
$wmain = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that used to sit right after the
#    kipalog.com hyperlink paragraph.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# Helper: append a "clean" paragraph (no rsid/style inheritance noise) built
# from a literal OOXML fragment at the very end of the document body.
# ---------------------------------------------------------------------------
function Append-Xml([string]$fragment) {
    $rng = $d.Content
    $rng.Collapse(0)
    $rng.InsertXML($fragment)
}

# Helper: append an empty paragraph, then turn it into a hyperlink paragraph
# (using Hyperlinks.Add so Word wires up the relationship + Hyperlink style
# the same way the real object model would).
function Append-HyperlinkParagraph([string]$url) {
    Append-Xml("<w:p $wmain/>")
    $endPos = $d.Content.End - 1
    $target = $d.Range($endPos, $endPos)
    $d.Hyperlinks.Add($target, $url, $null, $null, $url) | Out-Null
}

# ---------------------------------------------------------------------------
# 2) Append the new "Lo trinh hoc tieng Nhat" block at the end of the body.
# ---------------------------------------------------------------------------

# Separator line
Append-Xml("<w:p $wmain><w:r><w:t xml:space=`"preserve`">------------------------------------------------------------------------------------ </w:t></w:r></w:p>")

# Heading: "Con duong su nghiep cho Developer:"
Append-Xml("<w:p $wmain><w:pPr><w:rPr><w:b/><w:color w:val=`"00B050`"/><w:sz w:val=`"28`"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:color w:val=`"00B050`"/><w:sz w:val=`"28`"/></w:rPr><w:t>Con đường sự nghiệp cho Developer:</w:t></w:r></w:p>")

# Three career-path links
Append-HyperlinkParagraph("https://itviec.com/blog/con-duong-su-nghiep/")
Append-HyperlinkParagraph("https://itviec.com/blog/ky-su-cau-noi-la-gi/?utm_source=blogpost&utm_medium=referral&utm_content=con-duong-su-nghiep&utm_campaign=con-duong-su-nghiep-it")
Append-HyperlinkParagraph("https://viblo.asia/p/brse-co-can-thiet-phai-master-mot-ngon-ngu-lap-trinh-V3m5WBwWlO7")

# "Tu vung tieng Nhat IT:" paragraph - this is where the _GoBack bookmark
# moved to.
Append-Xml("<w:p $wmain><w:r><w:t>Từ vựng tiếng Nhật IT:</w:t></w:r><w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/></w:p>")

# Two vocabulary links
Append-HyperlinkParagraph("http://tiengnhat.minder.vn/tu-vung-tieng-nhat-chuyen-nganh-cong-nghe-thong-tin/")
Append-HyperlinkParagraph("https://kosaido-hr.com/articles/411")

# Three trailing empty paragraphs
Append-Xml("<w:p $wmain/>")
Append-Xml("<w:p $wmain/>")
Append-Xml("<w:p $wmain/>")

Write-Output "done"
